$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting rows 210-213 down to 211-214.
$ws.Rows.Item(210).Insert()

# Copy the style of the date cell from the (now shifted) row below so the
# new row's date cell keeps the same date/time number format.
$ws.Range("D211").Copy()
$ws.Range("D210").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values (same structure as surrounding rows).
$ws.Cells.Item(210, 1).Value = 4
$ws.Cells.Item(210, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(210, 3).Value = "Los Lagos"
$ws.Cells.Item(210, 4).Value = 44595
$ws.Cells.Item(210, 5).Value = 10
$ws.Cells.Item(210, 6).Value = 100112003
$ws.Cells.Item(210, 7).Value = "Ajo"
$ws.Cells.Item(210, 8).Value = "Chino"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 70
$ws.Cells.Item(210, 11).Value = 21000
$ws.Cells.Item(210, 12).Value = 21000
$ws.Cells.Item(210, 13).Value = 21000
$ws.Cells.Item(210, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(210, 15).Value = "China"
$ws.Cells.Item(210, 16).Value = 2100
$ws.Cells.Item(210, 17).Value = 10
$ws.Cells.Item(210, 18).Value = "Hortaliza"
